$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Applicant owns land" / "Permission obtained" / "Permission not obtained
# details" rows (65-67) are removed entirely; everything below shifts up by
# three rows. Deleting the rows (rather than clearing them) also keeps the
# merged-cell ranges and row numbering in sync automatically.
$ws.Rows("65:67").Delete()

# Column C narrows from 33 characters to 31 characters. Excel's ColumnWidth
# property is offset from the stored OOXML width by ~5/6 of a character, so
# compensate to land exactly on the target stored width of 31.
$ws.Columns("C").ColumnWidth = 31 - (5/6)
